# "file chooser implemented for the excel"
# - Fix trailing-space typo in B2 ("good " -> "good") and make it bold.
# - Add a new row (asd / no) below the existing data.
# - Move the active selection onto the new bottom-right area (C3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up "good " -> "good" and emphasize it with bold.
$ws.Range("B2").Value = "good"
$ws.Range("B2").Font.Bold = $true

# New data row.
$ws.Range("A3").Value = "asd"
$ws.Range("B3").Value = "no"

# Keep the selection in sync with the newly-used range.
[void]$ws.Range("C3").Select()
